$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the "temperature" column number format (scientific, style s="1")
# before A2 gets overwritten, so the new A7:A10 temperature cells - and the
# re-purposed B2 cell - can reuse it.
$tempFormat = $ws.Range("A2").NumberFormat

# --- New header cells (D1, E1) ---
$ws.Range("D1").Value = "OI Ratio"
$ws.Range("E1").Value = "NII Ratio"

# --- Row 2: temperature 10000, new B/C formulas, new D/E columns ---
$ws.Range("A2").Value = 10000
$ws.Range("B2").NumberFormat = $tempFormat
$ws.Range("B2").Formula = "=LOG10(15.7)"
$ws.Range("C2").Formula = "=LOG10(0.159/2.86)"
$ws.Range("D2").Formula = "=LOG10(0.0355/2.86)"
$ws.Range("E2").Formula = "=LOG10(1.01/2.86)"

# --- Row 3: temperature 100000 (was 420000), new formulas ---
$ws.Range("A3").Value = 100000
$ws.Range("B3").Formula = "=LOG10(12.9)"
$ws.Range("C3").Formula = "=LOG10(0.131/2.88)"
$ws.Range("D3").Formula = "=LOG10(0.0333/2.88)"
$ws.Range("E3").Formula = "=LOG10(0.962/2.88)"

# --- Row 4: temperature 1000000, new formulas ---
$ws.Range("A4").Value = 1000000
$ws.Range("B4").Formula = "=LOG10(14.7)"
$ws.Range("C4").Formula = "=LOG10(0.149/2.87)"
$ws.Range("D4").Formula = "=LOG10(0.0338/2.87)"
$ws.Range("E4").Formula = "=LOG10(0.979/2.87)"

# --- Row 5: temperature 10000000, new formulas ---
$ws.Range("A5").Value = 10000000
$ws.Range("B5").Formula = "=LOG10(21.4)"
$ws.Range("C5").Formula = "=LOG10(0.0904/2.86)"
$ws.Range("D5").Formula = "=LOG10(0.0488/2.86)"
$ws.Range("E5").Formula = "=LOG10(1.22/2.86)"

# --- Row 7 (new block): temperature 10000 ---
$ws.Range("A7").NumberFormat = $tempFormat
$ws.Range("A7").Value = 10000
$ws.Range("B7").Formula = "=LOG10(17.2)"
$ws.Range("C7").Formula = "=LOG10(0.159/3.07)"
$ws.Range("D7").Formula = "=LOG10((0.0512)/3.07)"
$ws.Range("E7").Formula = "=LOG10(1.33/3.07)"

# --- Row 8: temperature 100000 ---
$ws.Range("A8").NumberFormat = $tempFormat
$ws.Range("A8").Value = 100000
$ws.Range("B8").Formula = "=LOG10(13.9)"
$ws.Range("C8").Formula = "=LOG10(0.153/3.07)"
$ws.Range("D8").Formula = "=LOG10((0.045)/3.07)"
$ws.Range("E8").Formula = "=LOG10(1.22/3.07)"

# --- Row 9: temperature 1000000 ---
$ws.Range("A9").NumberFormat = $tempFormat
$ws.Range("A9").Value = 1000000
$ws.Range("B9").Formula = "=LOG10(14.2)"
$ws.Range("C9").Formula = "=LOG10(0.18/3.07)"
$ws.Range("D9").Formula = "=LOG10((0.0484)/3.07)"
$ws.Range("E9").Formula = "=LOG10(1.28/3.07)"

# --- Row 10: temperature 10000000 ---
$ws.Range("A10").NumberFormat = $tempFormat
$ws.Range("A10").Value = 10000000
$ws.Range("B10").Formula = "=LOG10(24.3)"
$ws.Range("C10").Formula = "=LOG10(0.109/3.07)"
$ws.Range("D10").Formula = "=LOG10((0.0798)/3.07)"
$ws.Range("E10").Formula = "=LOG10(1.76/3.07)"

# --- Selection / view state to match target ---
$ws.Range("C7").Select()

Write-Output "applied edits"
